$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.308.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4366"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.318"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.832.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.01%  "
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06526"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.26%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.271"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.342.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.058"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.037.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.304"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.220"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.965"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02357"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.207"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.167"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.195"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.435"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.027"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07015"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.70%  "
